$d = $word.ActiveDocument
$br = [char]11

# 1. Insert the date line before "From: Austing Dong" (adds a new leading
#    segment plus two line breaks, matching the "+<w:t>May 22nd, 2023</w:t>..." hunk).
$r = $d.Range(0, 0)
$r.InsertBefore("May 22nd, 2023" + $br + $br)

# 2. Update the opening paragraph: "requirement." -> "requirements. I would
#    like to highlight the following for your consideration:"
$old = "I am writing to express my strong interest in applying for the position of IT Co-op Student at Windsor-Detroit Bridge Authority - Divisional Office. As a University of Waterloo Computer Science undergraduate student, I strongly believe that my technical competencies and academic background are closely in line with the job requirement."
$new = "I am writing to express my strong interest in applying for the position of IT Co-op Student at Windsor-Detroit Bridge Authority - Divisional Office. As a University of Waterloo Computer Science undergraduate student, I strongly believe that my technical competencies and academic background are closely in line with the job requirements. I would like to highlight the following for your consideration:"
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 3. Rewrite the "I loved computer science..." paragraph.
$old = "I loved computer science as well as developing applications since Middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest, and I did an excellent job in relevant courses in the beginning of my University studies. I found solving business challenges through programming is fascinating because this is the way I feel the sense of accomplishment. Such deep interest in programming and technology has motivated me to deep dive in related fields such as software development, quality assurance and machine learning."
$new = "My passion for computer science and application development began in middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest. This enthusiasm continued into my university studies, where I excelled in relevant courses. I find solving business challenges through programming fascinating, as it provides me with a sense of accomplishment. This deep interest has motivated me to explore related fields such as software development, quality assurance, and machine learning."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 4. Split the big "co-op program" paragraph into three paragraphs: the
#    co-op/GitHub paragraph, the project description paragraph, and the
#    "through understanding requirements..." paragraph (with its tail
#    merged into the confidence sentence).
$old = "The computer science co-op program at the University of Waterloo offered me a unique opportunity to take on both programming and logical courses. Through working on a massive number of technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including but not limited to object-oriented programming, web application development, artificial intelligence, algorithm design and data abstraction. Such projects can be viewed on my GitHub: https://github.com/AustingDong. One of the biggest projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information which helps users quickly locate their desired items based on keywords' weight. This application can be used to quickly get all the important items and keywords from NASA Technical Report Server which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project detail can be found here. Through understanding the project requirements, researching on coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills through interpreting and explaining technical concepts to my teammates while working in a team environment."
$new = "The computer science co-op program at the University of Waterloo has offered me a unique opportunity to take on both programming and logical courses. Through working on numerous technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including object-oriented programming, web application development, artificial intelligence, algorithm design, and data abstraction. My projects can be viewed on my GitHub: https://github.com/AustingDong." + $br + $br + "One of the most significant projects I led and built was an application that uses AI to extract keywords from articles containing scientific or technical information. This application helps users quickly locate their desired items based on keyword weight and can be used to efficiently retrieve important items and keywords from NASA Technical Report Server, which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project details can be found here." + $br + $br + "Through understanding project requirements, researching coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills, which I have demonstrated by interpreting and explaining technical concepts to my teammates while working in a team environment. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 5. Merge the old "I am extremely interested..." confidence sentence (now
#    relocated into the previous paragraph) with the "I am willing to
#    answer..." paragraph into a single paragraph.
$old = "I am extremely interested in advancing my career and contributing my skills to Windsor-Detroit Bridge Authority - Divisional Office. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment."
$new = "I am extremely interested in advancing my career and contributing my skills to Windsor-Detroit Bridge Authority - Divisional Office. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any questions or require additional information."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 6. Remove the now-obsolete standalone "I am willing to answer..." paragraph
#    (it was merged into the paragraph above) along with its leading break pair.
#    Find.Execute doesn't reliably match a run of two consecutive manual line
#    breaks as literal search text, so locate the sentence first and then
#    widen the deletion range backwards over the two preceding w:br marks.
$rng = $d.Content
$rng.Find.Execute("I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any question or require additional information.") | Out-Null
$delRng = $d.Range($rng.Start - 2, $rng.End)
$delRng.Delete()
